$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with result/profit data
$ws.Range("G103").Value = "Acierto"
$ws.Range("H103").Value = 1.75

$ws.Range("G108").Value = "Fallo"
$ws.Range("H108").Value = -1

$ws.Range("G118").Value = "Acierto"
$ws.Range("H118").Value = 2.5

# Append new row 125
$ws.Range("A125").Value = 14552907
$ws.Range("B125").Value = "'2025-09-04"
$ws.Range("B125").Style = "Normal"
$ws.Range("C125").Value = "Antoine Escoffier"
$ws.Range("D125").Value = "Mae Malige"
$ws.Range("E125").Value = "Gana Mae Malige"
$ws.Range("F125").Value = 2.75

# Append new row 126
$ws.Range("A126").Value = 14560147
$ws.Range("B126").Value = "'2025-09-04"
$ws.Range("B126").Style = "Normal"
$ws.Range("C126").Value = "Lanlana Tararudee"
$ws.Range("D126").Value = "Veronika Erjavec"
$ws.Range("E126").Value = "Gana Lanlana Tararudee"
$ws.Range("F126").Value = 2.5
